$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Animate the "ajouter/modifier/masquer ingrédient" rows (31-33) -> mark "style" (F) done
$ws.Range("F31").Value = $true
$ws.Range("F32").Value = $true
$ws.Range("F33").Value = $true

# "liste des catégories" row (34) -> un-check pseudo code / code / lien bd
$ws.Range("C34").Value = $false
$ws.Range("D34").Value = $false
$ws.Range("E34").Value = $false

# "ajouter catégorie" row (35) -> mark "style" (F) done
$ws.Range("F35").Value = $true

# Update the active selection to match the latest edited cell
$ws.Range("F35").Select()
